# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '97.418.69'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '3.723.42'
$ws.Range('E3').Value = '  +1.39%  '

$ws.Range('E4').Value = '  +0.01%  '

Set-TextValue $ws.Range('D5') '237.13'
$ws.Range('E5').Value = '  -0.94%  '

Set-TextValue $ws.Range('D6') '1.94'
$ws.Range('E6').Value = '  +3.49%  '

Set-TextValue $ws.Range('D7') '657.76'
$ws.Range('E7').Value = '  +0.44%  '

Set-TextValue $ws.Range('D8') '0.436'
$ws.Range('E8').Value = '  +3.20%  '

$ws.Range('E9').Value = '  -1.31%  '

$ws.Range('D11').Value = '3.722.91'
$ws.Range('E11').Value = '  +1.41%  '

$ws.Range('E12').Value = '  +19.07%  '

Set-TextValue $ws.Range('D13') '44.84'
$ws.Range('E13').Value = '  -1.21%  '

$ws.Range('E14').Value = '  +0.94%  '

Set-TextValue $ws.Range('D15') '6.91'
$ws.Range('E15').Value = '  +1.34%  '

$ws.Range('D16').Value = '4.416.79'
$ws.Range('E16').Value = '  +1.37%  '

$ws.Range('D17').Value = '97.083.06'
$ws.Range('E17').Value = '  +0.52%  '

Set-TextValue $ws.Range('D18') '8.99'
$ws.Range('E18').Value = '  +0.90%  '

$ws.Range('D19').Value = '3.722.38'
$ws.Range('E19').Value = '  +1.18%  '

Set-TextValue $ws.Range('D20') '13.07'
$ws.Range('E20').Value = '  +2.20%  '

Set-TextValue $ws.Range('D21') '18.73'
$ws.Range('E21').Value = '  -0.63%  '

Set-TextValue $ws.Range('D22') '0.507'
$ws.Range('E22').Value = '  -4.04%  '

Set-TextValue $ws.Range('D23') '525.58'
$ws.Range('E23').Value = '  -1.62%  '

$ws.Range('E24').Value = '  -0.48%  '

Set-TextValue $ws.Range('D25') '0.0000227'
$ws.Range('E25').Value = '  +11.47%  '

Set-TextValue $ws.Range('D26') '6.90'
$ws.Range('E26').Value = '  -4.46%  '

Set-TextValue $ws.Range('D27') '106.16'
$ws.Range('E27').Value = '  +3.57%  '

Set-TextValue $ws.Range('D28') '0.193'
$ws.Range('E28').Value = '  +16.15%  '

$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.920.51'
$ws.Range('E29').Value = '  +1.37%  '

$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D30') '13.55'
$ws.Range('E30').Value = '  -0.04%  '

Set-TextValue $ws.Range('D31') '12.63'
$ws.Range('E31').Value = '  +2.01%  '

Set-TextValue $ws.Range('D32') '3.02'
$ws.Range('E32').Value = '  -0.49%  '

$ws.Range('E33').Value = '  +0.07%  '

$ws.Range('E34').Value = '  +3.84%  '

Set-TextValue $ws.Range('D35') '1.83'
$ws.Range('E35').Value = '  -3.38%  '

Set-TextValue $ws.Range('D36') '32.62'
$ws.Range('E36').Value = '  -0.23%  '

Set-TextValue $ws.Range('D37') '0.998'
$ws.Range('E37').Value = '  -0.50%  '

Set-TextValue $ws.Range('D38') '642.33'
$ws.Range('E38').Value = '  -2.42%  '

$ws.Range('E39').Value = '  -0.71%  '

Set-TextValue $ws.Range('D40') '8.75'
$ws.Range('E40').Value = '  -1.72%  '

$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('E42').Value = '  +3.16%  '

Set-TextValue $ws.Range('D43') '0.499'
$ws.Range('E43').Value = '  +12.23%  '

Set-TextValue $ws.Range('D44') '6.77'
$ws.Range('E44').Value = '  +0.40%  '

Set-TextValue $ws.Range('D45') '40.69'
$ws.Range('E45').Value = '  +4.99%  '

$ws.Range('E46').Value = '  +1.85%  '

Set-TextValue $ws.Range('D47') '0.969'
$ws.Range('E47').Value = '  +0.55%  '

Set-TextValue $ws.Range('D48') '0.0457'
$ws.Range('E48').Value = '  -0.25%  '

Set-TextValue $ws.Range('D49') '2.40'
$ws.Range('E49').Value = '  +2.89%  '

Set-TextValue $ws.Range('D50') '23.64'
$ws.Range('E50').Value = '  +0.06%  '

Set-TextValue $ws.Range('D51') '8.65'
$ws.Range('E51').Value = '  -0.75%  '
